$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compañias")

# Rename the "...Id" template placeholders to their non-Id counterparts.
# Order matters: it controls the order new shared strings are appended in,
# which needs to line up with the target workbook (PrecioLista, Promociones,
# Procedencia).
$ws.Range("B13").Value = "{{Compañias.PrecioLista}}"
$ws.Range("B15").Value = "{{Compañias.Promociones}}"
$ws.Range("B11").Value = "{{Compañias.Procedencia}}"

# Move the sheet selection from B13:C13 to B11:C11.
[void]$ws.Range("B11:C11").Select()
